$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): values + single combined style ----
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Pases"
$ws.Range("C1").Value = "Frase"
$ws.Range("D1").Value = "linkurl"
$ws.Range("E1").Value = "linkinvitacion"
$ws.Range("F1").Value = "Link de Google Site"
$ws.Range("G1").Value = "Link de Invitacion"

$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---- Data row (row 2): values (text-typed, matching source export) ----
$ws.Range("A2").Value = "Familia González Rendón"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "7"

$ws.Range("C2").Value = "Son indispensables en éste día!"

# Empty-string cells (kept present, not removed) - quote-prefix trick
# yields a real empty Text value without deleting the cell.
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"

$ws.Range("F2").Value = "URL del HTML"
$ws.Range("G2").Value = "LINK DEL HTML"
